$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "18/03/2023"
$ws.Range("D3").Value = 183.3
$ws.Range("D4").Value = 191.1
$ws.Range("D5").Value = 219.7
$ws.Range("D6").Value = 235.3
$ws.Range("D7").Value = 161.2
$ws.Range("D8").Value = 184.6
$ws.Range("D9").Value = 143
$ws.Range("D10").Value = 135.2
$ws.Range("D11").Value = 139.1
$ws.Range("D12").Value = 130
$ws.Range("D13").Value = 78
$ws.Range("D14").Value = 42.9
$ws.Range("D18").Value = 3
$ws.Range("D19").Value = 11
$ws.Range("D20").Value = 28
$ws.Range("D21").Value = 18
$ws.Range("D22").Value = 25
$ws.Range("D23").Value = 14
$ws.Range("D24").Value = 18
$ws.Range("D25").Value = 10
$ws.Range("D26").Value = 12
$ws.Range("D27").Value = 10
$ws.Range("D28").Value = 12
$ws.Range("D29").Value = 5
$ws.Range("D30").Value = 4
$ws.Range("D34").Value = 26
$ws.Range("D35").Value = 43
$ws.Range("D36").Value = 73
$ws.Range("D37").Value = 153
$ws.Range("D38").Value = 128
$ws.Range("D39").Value = 117
$ws.Range("D40").Value = 109
$ws.Range("D41").Value = 93
$ws.Range("D42").Value = 91
$ws.Range("D43").Value = 118
$ws.Range("D44").Value = 97
$ws.Range("D45").Value = 99
$ws.Range("D46").Value = 85
$ws.Range("D47").Value = 76
$ws.Range("D48").Value = 61
$ws.Range("D49").Value = 24
$ws.Range("D50").Value = 24
$ws.Range("D51").Value = 19
$ws.Range("D52").Value = 55
$ws.Range("D53").Value = 80
$ws.Range("D54").Value = 86
$ws.Range("D55").Value = 93
$ws.Range("D56").Value = 96
$ws.Range("D57").Value = 72
$ws.Range("D58").Value = 72
$ws.Range("D59").Value = 57
$ws.Range("D60").Value = 64
$ws.Range("D61").Value = 59
$ws.Range("D62").Value = 55
$ws.Range("D63").Value = 37
$ws.Range("D68").Value = 18
$ws.Range("D69").Value = 27
$ws.Range("D70").Value = 28
$ws.Range("D71").Value = 33
$ws.Range("D72").Value = 32
$ws.Range("D73").Value = 29
$ws.Range("D74").Value = 23
$ws.Range("D75").Value = 24
$ws.Range("D76").Value = 26
$ws.Range("D77").Value = 24
$ws.Range("D78").Value = 20
$ws.Range("D79").Value = 14
$ws.Range("D80").Value = 8
$ws.Range("D81").Value = 5
$ws.Range("D82").Value = 3
$ws.Range("D84").Value = 9
$ws.Range("D87").Value = 2
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 2
$ws.Range("D90").Value = 11
$ws.Range("D91").Value = 61
$ws.Range("D92").Value = 67
$ws.Range("D93").Value = 102
$ws.Range("D94").Value = 136
$ws.Range("D95").Value = 143
$ws.Range("D96").Value = 147
$ws.Range("D97").Value = 116
$ws.Range("D98").Value = 124
$ws.Range("D99").Value = 110
$ws.Range("D100").Value = 107
$ws.Range("D101").Value = 101
$ws.Range("D102").Value = 100
$ws.Range("D103").Value = 76
$ws.Range("D104").Value = 42
$ws.Range("D105").Value = 32
$ws.Range("D106").Value = 16
$ws.Range("D107").Value = 8
$ws.Range("D109").Value = 12
$ws.Range("D110").Value = 25
$ws.Range("D111").Value = 20
$ws.Range("D112").Value = 18
$ws.Range("D113").Value = 11
$ws.Range("D114").Value = 9
$ws.Range("D115").Value = 7
$ws.Range("D116").Value = 8
$ws.Range("D117").Value = 7
$ws.Range("D118").Value = 5
$ws.Range("D119").Value = 6
$ws.Range("D120").Value = 2
$ws.Range("D125").Value = 0
$ws.Range("D128").Value = 0
$ws.Range("D129").Value = 0
$ws.Range("D131").Value = 0
$ws.Range("D134").Value = 0
$ws.Range("D135").Value = 0
$ws.Range("D140").Value = 4
$ws.Range("D141").Value = 6
$ws.Range("D142").Value = 8
$ws.Range("D143").Value = 9
$ws.Range("D144").Value = 10
$ws.Range("D145").Value = 8
$ws.Range("D146").Value = 7
$ws.Range("D147").Value = 8
$ws.Range("D148").Value = 6
$ws.Range("D149").Value = 8
$ws.Range("D150").Value = 9
$ws.Range("D151").Value = 4
$ws.Range("D152").Value = 2
$ws.Range("D153").Value = 1
$ws.Range("D154").Value = 1
